# Auto-generated edit script applying the Midgardsormr_Profits price-refresh diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the active workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7816.5884
$ws.Range("I62").Value = 6918.5
$ws.Range("K62").Value = 6918.5
$ws.Range("M62").Value = -6294.5

$ws.Range("H65").Value = 7816.5884
$ws.Range("I65").Value = 6918.5
$ws.Range("K65").Value = 34592.5
$ws.Range("M65").Value = -31472.5

$ws.Range("H96").Value = 8929636
$ws.Range("I96").Value = 23809748
$ws.Range("K96").Value = 71429244
$ws.Range("M96").Value = -71427871

$ws.Range("H100").Value = 70022.336
$ws.Range("I100").Value = 77537.625
$ws.Range("J100").Value = 9900
$ws.Range("K100").Value = 77537.625
$ws.Range("L100").Value = 9900
$ws.Range("M100").Value = -76996.625
$ws.Range("N100").Value = -10982

$ws.Range("H107").Value = 4005.3076
$ws.Range("J107").Value = 7184.25
$ws.Range("L107").Value = 7184.25
$ws.Range("N107").Value = -11024.25

$ws.Range("H132").Value = 5983722
$ws.Range("I132").Value = 6243797
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 18731391
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -18728861
$ws.Range("N132").Value = -11060

$ws.Range("H137").Value = 8846.187
$ws.Range("I137").Value = 13666.275
$ws.Range("J137").Value = 4186.7666
$ws.Range("K137").Value = 40998.825
$ws.Range("L137").Value = 12560.2998
$ws.Range("M137").Value = -38448.825
$ws.Range("N137").Value = -17660.2998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5542.5
$ws.Range("I32").Value = 5271.5054
$ws.Range("K32").Value = 5271.5054
$ws.Range("M32").Value = -4984.5054

$ws.Range("H61").Value = 2864.8704
$ws.Range("I61").Value = 2038.2667
$ws.Range("J61").Value = 6997.8887
$ws.Range("K61").Value = 2038.2667
$ws.Range("L61").Value = 6997.8887
$ws.Range("M61").Value = -1826.2667
$ws.Range("N61").Value = -7421.8887

$ws.Range("H97").Value = 2803.1025
$ws.Range("I97").Value = 2430.6553
$ws.Range("K97").Value = 2430.6553
$ws.Range("M97").Value = -1934.6553

$ws.Range("H102").Value = 4016.5789
$ws.Range("I102").Value = 3132.2307
$ws.Range("K102").Value = 3132.2307
$ws.Range("M102").Value = -1510.2307

$ws.Range("H122").Value = 2233.2812
$ws.Range("I122").Value = 2052.4285
$ws.Range("K122").Value = 6157.2855
$ws.Range("M122").Value = -3707.2855

$ws.Range("H136").Value = 2864.8704
$ws.Range("I136").Value = 2038.2667
$ws.Range("J136").Value = 6997.8887
$ws.Range("K136").Value = 6114.800099999999
$ws.Range("L136").Value = 20993.6661
$ws.Range("M136").Value = -3564.800099999999
$ws.Range("N136").Value = -26093.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15029.921
$ws.Range("I20").Value = 22622.125
$ws.Range("J20").Value = 2014.7142
$ws.Range("K20").Value = 22622.125
$ws.Range("L20").Value = 2014.7142
$ws.Range("M20").Value = -22375.125
$ws.Range("N20").Value = -2508.7142

$ws.Range("H75").Value = 35000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 35000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H86").Value = 11726
$ws.Range("I86").Value = 1452.6
$ws.Range("K86").Value = 1452.6
$ws.Range("M86").Value = -329.5999999999999

$ws.Range("H89").Value = 11726
$ws.Range("I89").Value = 1452.6
$ws.Range("K89").Value = 7263
$ws.Range("M89").Value = -1647

$ws.Range("H94").Value = 1926.6957
$ws.Range("I94").Value = 1979.4286
$ws.Range("J94").Value = 1373
$ws.Range("K94").Value = 1979.4286
$ws.Range("L94").Value = 1373
$ws.Range("M94").Value = -1528.4286
$ws.Range("N94").Value = -2275

$ws.Range("H99").Value = 5453.1816
$ws.Range("I99").Value = 5776.1113
$ws.Range("K99").Value = 5776.1113
$ws.Range("M99").Value = -4278.1113

$ws.Range("H107").Value = 17636.818
$ws.Range("I107").Value = 20501.148
$ws.Range("K107").Value = 20501.148
$ws.Range("M107").Value = -18581.148

$ws.Range("H134").Value = 2412.3396
$ws.Range("I134").Value = 2290.745
$ws.Range("K134").Value = 6872.235
$ws.Range("M134").Value = -4337.235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 5449.5
$ws.Range("I21").Value = 9000
$ws.Range("J21").Value = 1899
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 1899
$ws.Range("M21").Value = -8765
$ws.Range("N21").Value = -2369

$ws.Range("H22").Value = 1037.25
$ws.Range("I22").Value = 1024.75
$ws.Range("J22").Value = 1049.75
$ws.Range("K22").Value = 1024.75
$ws.Range("L22").Value = 1049.75
$ws.Range("M22").Value = -674.75
$ws.Range("N22").Value = -1749.75

$ws.Range("H31").Value = 3653.75
$ws.Range("I31").Value = 3078.8708
$ws.Range("J31").Value = 4193.788
$ws.Range("K31").Value = 3078.8708
$ws.Range("L31").Value = 4193.788
$ws.Range("M31").Value = -2783.8708
$ws.Range("N31").Value = -4783.788

$ws.Range("H34").Value = 3653.75
$ws.Range("I34").Value = 3078.8708
$ws.Range("J34").Value = 4193.788
$ws.Range("K34").Value = 3078.8708
$ws.Range("L34").Value = 4193.788
$ws.Range("M34").Value = -2876.8708
$ws.Range("N34").Value = -4597.788

$ws.Range("H134").Value = 2160.1365
$ws.Range("I134").Value = 1572.9445
$ws.Range("K134").Value = 4718.833500000001
$ws.Range("M134").Value = -2183.833500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 382.6154
$ws.Range("J2").Value = 177.2
$ws.Range("L2").Value = 1063.2
$ws.Range("N2").Value = -1289.2

$ws.Range("H22").Value = 4910.6787
$ws.Range("J22").Value = 4999.963
$ws.Range("L22").Value = 14999.889
$ws.Range("N22").Value = -15337.889

$ws.Range("H27").Value = 4910.6787
$ws.Range("J27").Value = 4999.963
$ws.Range("L27").Value = 14999.889
$ws.Range("N27").Value = -15203.889

$ws.Range("H34").Value = 4268.923
$ws.Range("I34").Value = 886.75
$ws.Range("J34").Value = 4883.864
$ws.Range("K34").Value = 2660.25
$ws.Range("L34").Value = 14651.592
$ws.Range("M34").Value = -2576.25
$ws.Range("N34").Value = -14819.592

$ws.Range("H39").Value = 4097.6665
$ws.Range("I39").Value = 1144.5
$ws.Range("J39").Value = 10004
$ws.Range("K39").Value = 3433.5
$ws.Range("L39").Value = 30012
$ws.Range("M39").Value = -3139.5
$ws.Range("N39").Value = -30600

$ws.Range("H55").Value = 3998.8
$ws.Range("J55").Value = 4916.6665
$ws.Range("L55").Value = 14749.9995
$ws.Range("N55").Value = -15103.9995

$ws.Range("H105").Value = 9914.286
$ws.Range("J105").Value = 9914.286
$ws.Range("L105").Value = 29742.858
$ws.Range("N105").Value = -34984.858

$ws.Range("H132").Value = 1575.0312
$ws.Range("I132").Value = 1689.4
$ws.Range("J132").Value = 1474.1177
$ws.Range("K132").Value = 15204.6
$ws.Range("L132").Value = 13267.0593
$ws.Range("M132").Value = -12674.6
$ws.Range("N132").Value = -18327.0593

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 16000000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -171

$ws.Range("H132").Value = 1546.931
$ws.Range("I132").Value = 1588.25
$ws.Range("J132").Value = 390
$ws.Range("K132").Value = 4764.75
$ws.Range("L132").Value = 1170
$ws.Range("M132").Value = -2234.75
$ws.Range("N132").Value = -6230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1155.6207
$ws.Range("I22").Value = 760.6667
$ws.Range("K22").Value = 760.6667
$ws.Range("M22").Value = -465.6667

$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20686

$ws.Range("H25").Value = 18303.818
$ws.Range("I25").Value = 18935.2
$ws.Range("K25").Value = 18935.2
$ws.Range("M25").Value = -18705.2

$ws.Range("H27").Value = 1155.6207
$ws.Range("I27").Value = 760.6667
$ws.Range("K27").Value = 760.6667
$ws.Range("M27").Value = -653.6667

$ws.Range("H58").Value = 3999.75
$ws.Range("J58").Value = 3999.75
$ws.Range("L58").Value = 3999.75
$ws.Range("N58").Value = -4519.75

$ws.Range("H61").Value = 1468.3334
$ws.Range("I61").Value = 1465
$ws.Range("K61").Value = 1465
$ws.Range("M61").Value = -1263

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H113").Value = 1468.3334
$ws.Range("I113").Value = 1465
$ws.Range("K113").Value = 1465
$ws.Range("M113").Value = 705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 8642
$ws.Range("J22").Value = 17034
$ws.Range("L22").Value = 17034
$ws.Range("N22").Value = -17620

$ws.Range("H70").Value = 70009.125
$ws.Range("J70").Value = 77139.71000000001
$ws.Range("L70").Value = 77139.71000000001
$ws.Range("N70").Value = -77769.71000000001

$ws.Range("H73").Value = 70009.125
$ws.Range("J73").Value = 77139.71000000001
$ws.Range("L73").Value = 77139.71000000001
$ws.Range("N73").Value = -79323.71000000001

$ws.Range("H96").Value = 1654.1666
$ws.Range("I96").Value = 1899.8
$ws.Range("K96").Value = 1899.8
$ws.Range("M96").Value = -526.8

Write-Output "Updated cells: 241 set, 4 cleared"